$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PEIIR-IFPEIIR")

# Copy the format of the last existing header cell (K1) onto the two new
# header cells (L1, M1) before writing their text so they pick up the same
# style (right-aligned, wrapped, bold header font) instead of the sheet's
# plain default style.
$ws.Range("K1").Copy()
$ws.Range("L1:M1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("L1").Value = "green hydrogen"
$ws.Range("M1").Value = "low carbon hydrogen"

for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 12).Value = 0
    $ws.Cells.Item($r, 13).Value = 0
}

# New columns L (green hydrogen) and M (low carbon hydrogen) take the same
# column width as the preceding K column.
$ws.Columns.Item(12).ColumnWidth = $ws.Columns.Item(11).ColumnWidth
$ws.Columns.Item(13).ColumnWidth = $ws.Columns.Item(11).ColumnWidth

# Make the Industrial Fuel PEIIR sheet the active tab, with L2 selected,
# mirroring the saved workbook state after the edit.
$ws.Activate()
[void]$ws.Range("L2").Select()
